$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.383.19'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.46%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.915.27'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.19%  '

# Row 4
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.38%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.721'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +9.99%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '251.50'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.17%  '

# Row 7
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.34%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '40.69'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.01%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.358'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.35%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '53.02'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +7.24%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0733'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.99%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0999'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.27%  '

# Row 13
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.27%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.65'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.13%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.716'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.40%  '

# Row 16
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.925.01'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.02%  '

# Row 17
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.91'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.46%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '35.404.94'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.35%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.14'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.01%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0829'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.90%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.08'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +3.55%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '241.77'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.40%  '

# Row 23
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +4.74%  '

# Row 24
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.31%  '

# Row 25
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.32'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.61%  '

# Row 26
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.36'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +6.72%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '167.45'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.30%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.68'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.27%  '

# Row 29
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +5.18%  '

# Row 30
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.63%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.127.42'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +19.42%  '

# Row 32
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +4.43%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.97'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +12.44%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0580'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.81%  '

# Row 35
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'TrustWalletToken'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.58'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +18.12%  '

# Row 36
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'InternetComputer(DFINITY)'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.25'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.47%  '

# Row 37
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.35%  '

# Row 38
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.20%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.07'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.83%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.56'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +11.30%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.54'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +8.05%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.13'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.67%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0209'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.05%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0650'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.53%  '

# Row 45
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.75%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.346.48'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.51%  '

# Row 47
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.75%  '

# Row 48
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.51%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '45.54'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -4.93%  '

# Row 50
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.82%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.93'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -6.02%  '
